# Replace the original "Sheet1" with a fresh sheet named "Planilha1"
# (mirrors the author's workflow: new worksheet added, old one removed --
# which is also why Excel bumps sheetId from 1 to 2 in the saved file).
$wb = $excel.ActiveWorkbook
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Planilha1"
$oldSheet = $wb.Worksheets.Item("Sheet1")
$oldSheet.Delete()

$ws = $wb.Worksheets.Item("Planilha1")

# ---- A2 first (matches the original authoring order: the shared-string
# table in the target ends up with "A" at index 0, "Estrutura" at index 1,
# etc., which only happens if A2 is written before the header row). ----
$ws.Range("A2").Value = "A"

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "Estrutura"
$ws.Range("B1").Value = " Níveis"
$ws.Range("C1").Value = " Membros_Nivel_1"
$ws.Range("D1").Value = " Membros_Nivel_2"
$ws.Range("E1").Value = " Membros_Nivel_3"
$ws.Range("F1").Value = " Numero_de_Departamentos"
$ws.Range("G1").Value = " Membros"
$ws.Range("H1").Value = " Org_Size"
$ws.Range("I1").Value = " CC"

# ---- Remaining data row (row 2) values ----
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 6
$ws.Range("F2").Value = 2
$ws.Range("G2").Value = 10
$ws.Range("H2").Value = 100
$ws.Range("I2").Value = 200

# ---- Formatting ----
# Center every populated cell (both rows share the same centered alignment).
$used = $ws.Range("A1:I2")
$used.HorizontalAlignment = -4108
$used.VerticalAlignment = -4108

# Header row additionally gets a Text ("@") number format.
$ws.Range("A1:I1").NumberFormat = "@"

# ---- Column widths (approximate best-fit, matching the author's autosize) ----
$ws.Columns.Item(1).ColumnWidth = 49 / 6
$ws.Columns.Item(2).ColumnWidth = 4625 / 768
$ws.Columns.Item(3).ColumnWidth = 13073 / 768
$ws.Columns.Item(4).ColumnWidth = 13073 / 768
$ws.Columns.Item(5).ColumnWidth = 13073 / 768
$ws.Columns.Item(6).ColumnWidth = 20315 / 768
$ws.Columns.Item(7).ColumnWidth = 6929 / 768
$ws.Columns.Item(8).ColumnWidth = 1595 / 192
$ws.Columns.Item(9).ColumnWidth = 19 / 6

# ---- View state: selection on A2, page margins switched to the metric set ----
$ws.Range("A2").Select()

$ws.PageSetup.LeftMargin = $excel.CentimetersToPoints(1.3)
$ws.PageSetup.RightMargin = $excel.CentimetersToPoints(1.3)
$ws.PageSetup.TopMargin = $excel.CentimetersToPoints(2)
$ws.PageSetup.BottomMargin = $excel.CentimetersToPoints(2)
$ws.PageSetup.HeaderMargin = $excel.CentimetersToPoints(0.8)
$ws.PageSetup.FooterMargin = $excel.CentimetersToPoints(0.8)
